# 4.2.1.xlsx — capitalize the "by sex" / "mother's education" / "wealth
# quintile" sub-header rows (columns A-C) on the single worksheet, and
# reset the sheet's active cell back to A1 (it was left on A23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A/B first (ru/ky), matching the authoring order of the edit,
# then column C (en) afterwards.
$ws.Range("A14").Value = "Жынысы боюнча"
$ws.Range("B14").Value = "По полу"
$ws.Range("B17").Value = "Образование матери "
$ws.Range("A17").Value = "Энесинин билими "
$ws.Range("B23").Value = "Квинтиль по индексу благосостояния"
$ws.Range("C14").Value = "By sex"
$ws.Range("C17").Value = "Education of mother"
$ws.Range("C23").Value = "Wealth quintile"

# Restore the default selection to A1 (source file had A23 selected).
$ws.Range("A1").Select()
